# Update column A values for the specified rows (data imputation result update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = -21.60710000000001
$ws.Range("A32").Value = -21.26719999999999
$ws.Range("A36").Value = -20.2503
$ws.Range("A38").Value = -19.6711
$ws.Range("A46").Value = -21.7583
$ws.Range("A54").Value = -21.82999999999999
$ws.Range("A55").Value = -22.446
$ws.Range("A67").Value = -21.42809999999998
$ws.Range("A69").Value = -21.52749999999997
$ws.Range("A72").Value = -21.90219999999999
$ws.Range("A91").Value = -21.4265
$ws.Range("A99").Value = -20.19769999999999
